# Weekly price-sheet update: a new daily record is inserted above the
# existing row 87, pushing the remaining records (old rows 87-137) down
# by one row (new rows 88-138). The sheet's used range grows from
# A1:R137 to A1:R138 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 87..137 down to 88..138, leaving a blank row 87 in place.
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new weekly record.
$ws.Cells.Item(87, 1).Value = 4
$ws.Cells.Item(87, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(87, 3).Value = "Los Lagos"
$ws.Cells.Item(87, 4).Value = 44438
$ws.Cells.Item(87, 5).Value = 10
$ws.Cells.Item(87, 6).Value = 100112043
$ws.Cells.Item(87, 7).Value = "Pepino ensalada"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 150
$ws.Cells.Item(87, 11).Value = 20000
$ws.Cells.Item(87, 12).Value = 20000
$ws.Cells.Item(87, 13).Value = 20000
$ws.Cells.Item(87, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(87, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(87, 16).Value = 333
$ws.Cells.Item(87, 17).Value = 60
$ws.Cells.Item(87, 18).Value = "Hortaliza"
